$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34: new diary entry (02.12.18, Open External Software working) ---

# A34: date label, kept as literal text (matches the rest of column A)
$ws.Range("A34").NumberFormat = "@"
$ws.Range("A34").Value = "02.12.18"
$ws.Range("A34").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

# B34 / C34: begin / end times
$ws.Range("B34").Value = 0.45833333333333331
$ws.Range("C34").Value = 0.79166666666666663

# D34: elapsed time formula (same pattern as the rest of column D)
$ws.Range("D34").Formula = "=C34-B34"

# E34: remarks text, word-wrapped like the row above it
$ws.Range("E34").WrapText = $true
$ws.Range("E34").Value = "-Open External Software`n-Tooltipps"

$ws.Rows.Item(34).RowHeight = 30

# Update the active selection to reflect where the user ended up
$ws.Range("C35").Select()
